$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PCS Voltage value (C10) from 800 to 690
$ws.Range("C10").Value = 690

# Update the active cell selection to C11
[void]$ws.Range("C11").Select()
